$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for the 380 kV case (rows 2-25).
# Columns A, F, J, N are untouched (stay at their existing 0 / index values).
$newValues = @{
    2 = @{ "B" = 13.11113637047131; "C" = 5.705597563789397; "D" = 6.039425842009615; "E" = 11.13424970147505; "G" = 55.29294700333538; "H" = 20.72956375566558; "I" = 33.10481550870235; "K" = 11.21969019655359; "L" = 10.42753119365155; "M" = 14.8707300276895 }
    3 = @{ "B" = 12.9873261089986; "C" = 5.544432999783423; "D" = 5.930137760864997; "E" = 11.1488267138031; "G" = 55.01382656196402; "H" = 20.7233453029252; "I" = 33.0742167876713; "K" = 11.1307217502253; "L" = 10.43644797159606; "M" = 14.86966469510935 }
    4 = @{ "B" = 12.91466058696049; "C" = 5.441671711245212; "D" = 5.863844652023885; "E" = 11.1587426640961; "G" = 54.8523241640772; "H" = 20.72236886339234; "I" = 33.06008439908312; "K" = 11.07900631622095; "L" = 10.44337994173419; "M" = 14.87177433414062 }
    5 = @{ "B" = 12.88592716955062; "C" = 5.398866678726747; "D" = 5.837073109117473; "E" = 11.16302664164285; "G" = 54.78903489120108; "H" = 20.72268613536702; "I" = 33.05549801994846; "K" = 11.0586868137718; "L" = 10.4465714649012; "M" = 14.87333054304102 }
    6 = @{ "B" = 12.88121004218576; "C" = 5.391703695215847; "D" = 5.832643576690373; "E" = 11.16375268812781; "G" = 54.7786793176143; "H" = 20.72278202559843; "I" = 33.0548073308535; "K" = 11.0553590117491; "L" = 10.44712356961486; "M" = 14.87363105719145 }
    7 = @{ "B" = 12.91426947684265; "C" = 5.441098150142754; "D" = 5.863482563535138; "E" = 11.15879945437247; "G" = 54.85146035190687; "H" = 20.72237024595512; "I" = 33.06001779521512; "K" = 11.0787291944534; "L" = 10.4434214986689; "M" = 14.87179249975639 }
    8 = @{ "B" = 13.0677718252259; "C" = 5.650833870542392; "D" = 6.001598960925833; "E" = 11.13907564826676; "G" = 55.19468114436609; "H" = 20.72683027132177; "I" = 33.09330024992488; "K" = 11.18842371319091; "L" = 10.43030349375463; "M" = 14.86979009704127 }
    9 = @{ "B" = 13.39379162478419; "C" = 6.030829732240652; "D" = 6.277137670428607; "E" = 11.10804415744834; "G" = 55.94419449198028; "H" = 20.75808952705563; "I" = 33.1953951759151; "K" = 11.42556623340211; "L" = 10.41612487266826; "M" = 14.8877133871977 }
    10 = @{ "B" = 13.64630249940913; "C" = 6.28963542389709; "D" = 6.48015277394247; "E" = 11.0898871479247; "G" = 56.53867952798345; "H" = 20.79472221165235; "I" = 33.29267966414233; "K" = 11.61170961732543; "L" = 10.41272352209749; "M" = 14.91407903890372 }
    11 = @{ "B" = 13.7634960256498; "C" = 6.402718952735855; "D" = 6.572170986159237; "E" = 11.08263064390129; "G" = 56.8179820922637; "H" = 20.81433691814372; "I" = 33.34172394766041; "K" = 11.69863624963088; "L" = 10.41269321171981; "M" = 14.92890643075877 }
    12 = @{ "B" = 13.80816598451558; "C" = 6.444855797001588; "D" = 6.606931064382914; "E" = 11.08002668308346; "G" = 56.9249620807871; "H" = 20.82218663565276; "I" = 33.36097898055721; "K" = 11.73184691422614; "L" = 10.4128991793091; "M" = 14.93492529059152 }
    13 = @{ "B" = 13.79853323895842; "C" = 6.435811608612585; "D" = 6.599449222872439; "E" = 11.08058109728981; "G" = 56.90186902292701; "H" = 20.82047732572106; "I" = 33.35680178325682; "K" = 11.72468183384624; "L" = 10.41284516149727; "M" = 14.9336111090356 }
    14 = @{ "B" = 13.76716545834167; "C" = 6.406199403856803; "D" = 6.575032633301298; "E" = 11.08241353209442; "G" = 56.82675942472465; "H" = 20.81497427418732; "I" = 33.3432944164314; "K" = 11.70136280310008; "L" = 10.41270580313822; "M" = 14.92939352893859 }
    15 = @{ "B" = 13.74798844648999; "C" = 6.387971338425793; "D" = 6.560064595783902; "E" = 11.08355468366004; "G" = 56.78090889478834; "H" = 20.81165838978853; "I" = 33.33510953657582; "K" = 11.68711654141248; "L" = 10.41264873791515; "M" = 14.92686265160673 }
    16 = @{ "B" = 13.63868673746869; "C" = 6.282150113437726; "D" = 6.474129225591089; "E" = 11.09038153577322; "G" = 56.52059973475531; "H" = 20.79349953124381; "I" = 33.28957040584947; "K" = 11.60607146095989; "L" = 10.41275597100145; "M" = 14.91316676780922 }
    17 = @{ "B" = 13.57219861516537; "C" = 6.216028923228968; "D" = 6.421298063806791; "E" = 11.09482629516888; "G" = 56.36313719619915; "H" = 20.78311383691362; "I" = 33.26285641114756; "K" = 11.55690777396395; "L" = 10.41320982227445; "M" = 14.90548819911706 }
    18 = @{ "B" = 13.53417784342304; "C" = 6.177561113897704; "D" = 6.390881624519635; "E" = 11.09747725110651; "G" = 56.27340745578439; "H" = 20.77741810804603; "I" = 33.24794237604916; "K" = 11.52884364084207; "L" = 10.41361369772154; "M" = 14.90133861913286 }
    19 = @{ "B" = 13.52134400448938; "C" = 6.164462129239382; "D" = 6.380579258595081; "E" = 11.09839105098328; "G" = 56.24317235074092; "H" = 20.77553741489574; "I" = 33.24297038672429; "K" = 11.51937918649097; "L" = 10.41377499438992; "M" = 14.89997958657863 }
    20 = @{ "B" = 13.57925378268254; "C" = 6.223112957625876; "D" = 6.426925343456268; "E" = 11.09434337038596; "G" = 56.37981300063585; "H" = 20.78419067052447; "I" = 33.26565350813913; "K" = 11.56211944984355; "L" = 10.4131467312039; "M" = 14.906277993293 }
    21 = @{ "B" = 13.77637137288445; "C" = 6.414915954532878; "D" = 6.58220697787627; "E" = 11.08187139839299; "G" = 56.84878848056405; "H" = 20.81657921838388; "I" = 33.34724336834108; "K" = 11.70820444462964; "L" = 10.41274084060397; "M" = 14.9306213974826 }
    22 = @{ "B" = 13.90687717442152; "C" = 6.536268517174443; "D" = 6.68317978259913; "E" = 11.07455895376246; "G" = 57.16233826087428; "H" = 20.8402058307565; "I" = 33.40454521079928; "K" = 11.80537568218477; "L" = 10.41374262559337; "M" = 14.94888451102053 }
    23 = @{ "B" = 13.83708444550565; "C" = 6.471871655382712; "D" = 6.629347428785137; "E" = 11.07838511595617; "G" = 56.99436679128237; "H" = 20.82737168335392; "I" = 33.37360016075837; "K" = 11.75336834459269; "L" = 10.4130922690395; "M" = 14.93892302779276 }
    24 = @{ "B" = 13.5760635022951; "C" = 6.219911678108818; "D" = 6.42438138176397; "E" = 11.09456140286565; "G" = 56.37227138030328; "H" = 20.78370297673963; "I" = 33.26438755641945; "K" = 11.55976262364439; "L" = 10.41317480934773; "M" = 14.90592010184135 }
    25 = @{ "B" = 13.30314688562075; "C" = 5.931524215518912; "D" = 6.202334049541051; "E" = 11.11562237551427; "G" = 55.73352155635383; "H" = 20.74722773388851; "I" = 33.16384878480143; "K" = 11.35920679298945; "L" = 10.41872657393794; "M" = 14.88053687119638 }
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}

